$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin/Link swaps (rows 27/28, 49/50) ---
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

# --- Price (D) values that are NOT numeric-looking: set directly ---
$ws.Range('D2').Value = '60.631.36'
$ws.Range('D3').Value = '2.910.65'
$ws.Range('D9').Value = '2.911.73'
$ws.Range('D13').Value = '3.417.46'
$ws.Range('D15').Value = '60.775.26'
$ws.Range('D17').Value = '2.905.24'
$ws.Range('D26').Value = '3.039.18'
$ws.Range('D30').Value = '0.0₃0856'
$ws.Range('D41').Value = '2.333.98'

# --- Price (D) values that WOULD be auto-parsed as numbers: force text ---
$numericLookingRefs = @('D4', 'D5', 'D6', 'D8', 'D10', 'D11', 'D12', 'D16', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}
$ws.Range('D4').Value = '1.00'
$ws.Range('D5').Value = '527.67'
$ws.Range('D6').Value = '142.65'
$ws.Range('D8').Value = '0.552'
$ws.Range('D10').Value = '0.107'
$ws.Range('D11').Value = '5.85'
$ws.Range('D12').Value = '0.352'
$ws.Range('D16').Value = '22.61'
$ws.Range('D19').Value = '4.93'
$ws.Range('D20').Value = '11.52'
$ws.Range('D21').Value = '360.06'
$ws.Range('D22').Value = '6.45'
$ws.Range('D23').Value = '1.00'
$ws.Range('D24').Value = '5.66'
$ws.Range('D25').Value = '63.41'
$ws.Range('D27').Value = '0.449'
$ws.Range('D28').Value = '0.182'
$ws.Range('D31').Value = '7.62'
$ws.Range('D34').Value = '19.62'
$ws.Range('D35').Value = '152.93'
$ws.Range('D36').Value = '4.34'
$ws.Range('D37').Value = '5.58'
$ws.Range('D38').Value = '0.997'
$ws.Range('D42').Value = '1.46'
$ws.Range('D43').Value = '3.68'
$ws.Range('D44').Value = '0.643'
$ws.Range('D45').Value = '20.77'
$ws.Range('D46').Value = '0.0567'
$ws.Range('D47').Value = '1.00'
$ws.Range('D48').Value = '4.82'
$ws.Range('D49').Value = '10.37'
$ws.Range('D50').Value = '0.0233'
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).Style = "Normal"
}

# --- Volume(1h) (E) values: always text already (has % and spaces) ---
$ws.Range('E2').Value = '  -3.98%  '
$ws.Range('E3').Value = '  -3.73%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -5.67%  '
$ws.Range('E6').Value = '  -8.04%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('E9').Value = '  -3.89%  '
$ws.Range('E10').Value = '  -5.62%  '
$ws.Range('E11').Value = '  -9.17%  '
$ws.Range('E12').Value = '  -4.41%  '
$ws.Range('E13').Value = '  -3.76%  '
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('E15').Value = '  -3.84%  '
$ws.Range('E16').Value = '  -6.34%  '
$ws.Range('E17').Value = '  -4.13%  '
$ws.Range('E18').Value = '  -7.61%  '
$ws.Range('E19').Value = '  -3.65%  '
$ws.Range('E20').Value = '  -4.27%  '
$ws.Range('E21').Value = '  -9.64%  '
$ws.Range('E22').Value = '  -3.61%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('E26').Value = '  -3.78%  '
$ws.Range('E27').Value = '  -4.15%  '
$ws.Range('E28').Value = '  -4.26%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  -13.28%  '
$ws.Range('E31').Value = '  -12.57%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -4.25%  '
$ws.Range('E35').Value = '  -4.44%  '
$ws.Range('E36').Value = '  -8.66%  '
$ws.Range('E37').Value = '  -8.09%  '
$ws.Range('E38').Value = '  -10.26%  '
$ws.Range('E39').Value = '  -8.34%  '
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('E41').Value = '  -8.24%  '
$ws.Range('E42').Value = '  -7.97%  '
$ws.Range('E43').Value = '  -7.21%  '
$ws.Range('E44').Value = '  -4.12%  '
$ws.Range('E46').Value = '  -5.92%  '
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('E48').Value = '  -5.33%  '
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('E50').Value = '  -7.16%  '
$ws.Range('E51').Value = '  -2.88%  '
